# Work off the live ActiveSheet (going through $excel, not the detached $wb
# handle) so writes actually reach the model.
$ws = $excel.ActiveSheet

# Corrections to already-entered days
$ws.Range("C309").Value = 145
$ws.Range("C314").Value = 238
$ws.Range("L315").Value = 2
$ws.Range("C316").Value = 144
$ws.Range("C317").Value = 115

# Row 318 becomes a fully populated data row (new day's data)
$ws.Range("C318").Value = 21
$ws.Range("E318").Value = 13
$ws.Range("F318").Value = 8
$ws.Range("G318").Value = 82
$ws.Range("L318").Value = 0
$ws.Range("M318").Value = 0

# Update the active selection to reflect where the user ended up (bottom-right pane)
$ws.Range("Q17").Select()

$true
